# Article 35 ("King of my realm") is ready -- add it as a new row (row 36)
# to the articles table, and update the sheet view to focus on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ayat = @'
Surah Al Imran, 105 - 110
'@

$tags = @'
Imran Khan, Failing Economy, Calm Cities, Solution to Pakistan, Muslim Ummah, Unity and Harmony
'@

$content = @'
h1: King of my realm
p.note: I did 15 minutes yoga to organise these thoughts.. Listen please.
img.width-m-half: calm city.jpg
p.caption.text-center: Photo by <a href=https://unsplash.com/@nikolayv?utm_source=unsplash&utm_medium=referral&utm_content=creditCopyText target=_blank>Nikolay Vorobyev</a> on Unsplash
p: This city above has few high rise buildings across the ocean with low traffic rush. People are jogging in the evening. Some are going to meet their families. A lot of them are earning good. It speaks of comfort, luxury, peace and calmness. The effort of some good governances made it possible for these people to live in harmony and enjoy their lives. Good life expectancy, lower poverty, better education and above all; cleanliness. 
p: Good governance is directly linked with easier life. I am from Pakistan, so I take this country as my case study today. We are a progressing state. These past few years we have seen certain downs in this country. Economy is lowering, education system can not ignite passion in ourselves, poverty is increasing, each passing day it is becoming difficult to obtain calmness and peacefulness. Overall we all agree this country was nose diving into chaos. 
h3: What happened after the nose dive?
p: We got the best prime minister who have all traits which are needed to build this country back. Passionate, positive, not giving up, enlightened and believes in Pakistan transforming into Madina (a just state). Still he is talking heart winning speeches around the globe. Today he met one of the greatest leader, Mahatir Muhammad, and expressed himself maturely. Adequately delivered his message to all muslim states who were opposing the earlier summit held in Malaysia. 
img.width-m-half: imran statement in malaysia.png
p.caption.text-center: Link to article by Dawn is<a href=https://www.dawn.com/news/1532420/pm-regrets-missing-kl-summit-over-others-misconception-it-would-divide-ummah target=_blank> here</a>
h3: What does Quran say about ruling the realm?
quote: And do not be like the ones who became divided and differed after the clear proofs had come to them. And those will have a great punishment. <br> - Surah Al Imran verse 105
p: According to Quran, Imran Khan is working on the right path. His today’s statement emphasised upon unity of muslims, as are his previous efforts and statements. 
h3: Why should I care, I can’t find cheaper food.
p: Most among us are getting tired of future hopes and gradually we are getting buried in huge economic fall. Price of cars is going up, most commodities are getting harder to get our hands on and life is really not getting any easier. It is true, it is really not possible to live in peace with all these troubles. 
p: May be our Prime Minister needs to resolve internal issues before he gives lengthy talks on international problems. May be internal affairs take priority. May be if he places better people in his cabinet he is able to fix these issues sooner. 
p: All good things and correct decisions taken by Imran Khan will lead us to a better place. All good things and correct decisions taken by me do not really have an impact at national level. May be I should not really worry how we can get out of economic crises, because it is not my job. It is true, we really do not have to worry about economic crises in somebody’s elses realm. Realm of Imran Khan has to be fixed by himself. 
p: My personal opinion is, if Imran Khan keeps on talking about unity at international level, stays passionate and do not really give up, he will be able to address the poverty issue inside Pakistan also. It will take time.
img.width-m-half: Ghabrana nai hai.jpeg
h3: Coming back to the topic
p: I am king of my realm. My realm is my kitchen, my bathroom, my bedroom, entry lounge, lawn outside, the flowers in the pots, green garden, my part of the front road, neighbour, close relatives, colleagues, bosses, friends and people we interact everyday. I am king of all these areas. I need to talk of unity when there is a conflict that can be fixed with my intervention inside my realm. I sometimes do not have to talk, sometimes I have to act that brings all these pawns of my realm closer to each other. I do not need to say it out loud, ‘Pray your 5 prayers a day’. I need to just act that makes it colorful for the people of my realm to grasp it. 
p: I need to build the better building with mirrors glowing on top, people jogging in the evening, meeting their closed ones, where people feel financially comfortable and they really enjoy their lives. I need to fix the realm sooner, because people in my realm are suffering of conflicts and failures. 
p: If I start fixing my realm, people in my realm will create more kings. It will create a fusion factory where kings are made. As soon as we become a society of kings, we will see nearly no body who is poverty struck. We need to start from our branch of Imran Khan’s realm. Bottom up, slowly and gradually curing my realm, Pakistan, in’sha’Allah.
p: May Allah (swt) help us fix ourselves that we inspire our surroundings with our actions. Amen
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@

# Ser (row number within the table)
$ws.Range("A36").Value = 35

# Date published
$ws.Range("B36").Value = 43865
$ws.Range("B36").NumberFormat = "d-mmm-yy"

# Ayats (column C), Content (column D), Author (column E), Tags (column F)
# -- set in this order so the shared-string table is populated Ayats, Tags, Content
$ws.Range("C36").Value = $ayat
$ws.Range("F36").Value = $tags
$ws.Range("D36").Value = $content
$ws.Range("E36").Value = "Qasim Ali"

# Row 36 holds a very long wrapped article -- matches the auto-fit height used by
# the other long-article rows in this sheet.
$ws.Rows.Item(36).RowHeight = 409.6

# Scroll the view down to the new row and select its Content cell, like the author left it.
$ws.Range("D36").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 36

